$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "States"
$ws.Name = "States"

# Turn the A1:D52 range (Row, Column, State, Abbreviation) into an Excel Table
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:D52"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Move the active selection to C24
[void]$ws.Range("C24").Select()
